$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Log_Muestras")

$timestamps = @{
    2  = "2025-11-03T00:09:55.387251"
    3  = "2025-11-03T00:09:55.387251"
    4  = "2025-11-03T00:09:55.387251"
    5  = "2025-11-03T00:09:55.387251"
    6  = "2025-11-03T00:09:55.387251"
    7  = "2025-11-03T00:09:55.387251"
    8  = "2025-11-03T00:09:55.388252"
    9  = "2025-11-03T00:09:55.388252"
    10 = "2025-11-03T00:09:55.388252"
    11 = "2025-11-03T00:09:55.388252"
    12 = "2025-11-03T00:09:55.388252"
    13 = "2025-11-03T00:09:55.388252"
    14 = "2025-11-03T00:09:55.388252"
    15 = "2025-11-03T00:09:55.388252"
    16 = "2025-11-03T00:09:55.388252"
    17 = "2025-11-03T00:09:55.388252"
    18 = "2025-11-03T00:09:55.389252"
    19 = "2025-11-03T00:09:55.389252"
    20 = "2025-11-03T00:09:55.389252"
    21 = "2025-11-03T00:09:55.389252"
    22 = "2025-11-03T00:09:55.389252"
    23 = "2025-11-03T00:09:55.389252"
    24 = "2025-11-03T00:09:55.389252"
    25 = "2025-11-03T00:09:55.389252"
    26 = "2025-11-03T00:09:55.389252"
    27 = "2025-11-03T00:09:55.389252"
    28 = "2025-11-03T00:09:55.390251"
    29 = "2025-11-03T00:09:55.390251"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
